$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "Real Madrid  - Valencia CF: 19:00"
$ws.Cells.Item(2, 2).Value = "Real Madrid"
$ws.Cells.Item(2, 3).Value = 78
$ws.Cells.Item(2, 4).Value = 79
$ws.Cells.Item(2, 5).Value = 93
$ws.Cells.Item(2, 6).Value = 1.18
$ws.Cells.Item(2, 7).ClearContents()

# Row 3
$ws.Cells.Item(3, 1).Value = "Bayern Munich  - Bayer 04 Leverkusen: -:-'"
$ws.Cells.Item(3, 2).Value = "Bayern Munich"
$ws.Cells.Item(3, 3).Value = 74
$ws.Cells.Item(3, 4).Value = 81
$ws.Cells.Item(3, 5).Value = 85
$ws.Cells.Item(3, 6).Value = 1.23
$ws.Cells.Item(3, 7).ClearContents()

# Row 4
$ws.Cells.Item(4, 1).Value = "Club Brugge KV  - FCV Dender EH: -:-'"
$ws.Cells.Item(4, 2).Value = "Club Brugge KV"
$ws.Cells.Item(4, 3).Value = 73
$ws.Cells.Item(4, 4).Value = 86
$ws.Cells.Item(4, 5).Value = 100
$ws.Cells.Item(4, 6).Value = 1.28
$ws.Cells.Item(4, 7).ClearContents()

# Row 5
$ws.Cells.Item(5, 1).Value = "Galatasaray  - Trabzonspor: -:-'"
$ws.Cells.Item(5, 2).Value = "Galatasaray"
$ws.Cells.Item(5, 3).Value = 73
$ws.Cells.Item(5, 4).Value = 80
$ws.Cells.Item(5, 5).ClearContents()
$ws.Cells.Item(5, 6).Value = 1.44
$ws.Cells.Item(5, 7).ClearContents()

# Row 6
$ws.Cells.Item(6, 1).Value = "Haverfordwest County - The New Saints : 0:3'"
$ws.Cells.Item(6, 2).Value = "The New Saints"
$ws.Cells.Item(6, 3).Value = 71
$ws.Cells.Item(6, 4).Value = 76
$ws.Cells.Item(6, 5).Value = 85
$ws.Cells.Item(6, 6).Value = 1.61
$ws.Cells.Item(6, 7).ClearContents()

# Row 7
$ws.Cells.Item(7, 1).Value = "Olympiacos Piraeus  - Aris Thessaloniki: 17:00"
$ws.Cells.Item(7, 2).Value = "Olympiacos Piraeus"
$ws.Cells.Item(7, 3).Value = 69
$ws.Cells.Item(7, 4).Value = 78
$ws.Cells.Item(7, 5).ClearContents()
$ws.Cells.Item(7, 6).Value = 1.3
$ws.Cells.Item(7, 7).ClearContents()

# Row 8
$ws.Cells.Item(8, 1).Value = "Feyenoord Rotterdam  - FC Volendam: 18:00"
$ws.Cells.Item(8, 2).Value = "Feyenoord Rotterdam"
$ws.Cells.Item(8, 3).Value = 67
$ws.Cells.Item(8, 4).Value = 89
$ws.Cells.Item(8, 5).ClearContents()
$ws.Cells.Item(8, 6).Value = 1.12
$ws.Cells.Item(8, 7).ClearContents()

# Row 9
$ws.Cells.Item(9, 1).Value = "Club Africain  - AS Soliman: 1:1"
$ws.Cells.Item(9, 2).Value = "Club Africain"
$ws.Cells.Item(9, 3).Value = 63
$ws.Cells.Item(9, 4).Value = 86
$ws.Cells.Item(9, 5).Value = 100
$ws.Cells.Item(9, 6).Value = 1.28
$ws.Cells.Item(9, 7).ClearContents()

# Row 10
$ws.Cells.Item(10, 1).Value = "UD Almería  - SD Eibar: 19:00"
$ws.Cells.Item(10, 2).Value = "UD Almería"
$ws.Cells.Item(10, 3).Value = 63
$ws.Cells.Item(10, 4).ClearContents()
$ws.Cells.Item(10, 5).Value = 71
$ws.Cells.Item(10, 6).Value = 1.73
$ws.Cells.Item(10, 7).ClearContents()

# Row 11
$ws.Cells.Item(11, 1).Value = "Lee Man ✓ - Hong Kong Football Club: 7:1"
$ws.Cells.Item(11, 2).Value = "Lee Man"
$ws.Cells.Item(11, 3).Value = 60
$ws.Cells.Item(11, 4).Value = 75
$ws.Cells.Item(11, 5).ClearContents()
$ws.Cells.Item(11, 6).Value = 1.53
$ws.Cells.Item(11, 7).Value = "✓"

# Row 12
$ws.Cells.Item(12, 1).Value = "Al-Jazira Club  - Al-Bataeh CSC: -:-'"
$ws.Cells.Item(12, 2).Value = "Al-Jazira Club"
$ws.Cells.Item(12, 3).Value = 58
$ws.Cells.Item(12, 4).Value = 86
$ws.Cells.Item(12, 5).Value = 100
$ws.Cells.Item(12, 6).Value = 1.28
$ws.Cells.Item(12, 7).ClearContents()

# Row 13
$ws.Cells.Item(13, 1).Value = "Vitória Guimarães SC - SL Benfica : 19:30"
$ws.Cells.Item(13, 2).Value = "SL Benfica"
$ws.Cells.Item(13, 3).Value = 54
$ws.Cells.Item(13, 4).Value = 82
$ws.Cells.Item(13, 5).ClearContents()
$ws.Cells.Item(13, 6).Value = 1.5
$ws.Cells.Item(13, 7).ClearContents()

# Row 14
$ws.Cells.Item(14, 1).Value = "Atlético de Madrid  - Sevilla FC: 0:0'"
$ws.Cells.Item(14, 2).Value = "Atlético de Madrid"
$ws.Cells.Item(14, 3).Value = 54
$ws.Cells.Item(14, 4).Value = 79
$ws.Cells.Item(14, 5).Value = 93
$ws.Cells.Item(14, 6).Value = 1.18
$ws.Cells.Item(14, 7).ClearContents()

# Row 15
$ws.Cells.Item(15, 1).Value = "Liverpool FC Montevideo  - CA Juventud: 18:30"
$ws.Cells.Item(15, 2).Value = "Liverpool FC Montevideo"
$ws.Cells.Item(15, 3).Value = 54
$ws.Cells.Item(15, 4).Value = 51
$ws.Cells.Item(15, 5).Value = 67
$ws.Cells.Item(15, 6).Value = 1.69
$ws.Cells.Item(15, 7).ClearContents()

# Row 16
$ws.Cells.Item(16, 1).Value = "Tainan City TSG - Phnom Penh Crown ✓: 2:3"
$ws.Cells.Item(16, 2).Value = "Phnom Penh Crown"
$ws.Cells.Item(16, 3).Value = 53
$ws.Cells.Item(16, 4).Value = 70
$ws.Cells.Item(16, 5).ClearContents()
$ws.Cells.Item(16, 6).Value = 1.18
$ws.Cells.Item(16, 7).Value = "✓"

# Row 17
$ws.Cells.Item(17, 1).Value = "Club Alianza Lima  - FBC Melgar: 2:2"
$ws.Cells.Item(17, 2).Value = "Club Alianza Lima"
$ws.Cells.Item(17, 3).Value = 47
$ws.Cells.Item(17, 4).Value = 86
$ws.Cells.Item(17, 5).Value = 100
$ws.Cells.Item(17, 6).Value = 1.28
$ws.Cells.Item(17, 7).ClearContents()

# Row 18
$ws.Cells.Item(18, 1).Value = "Colwyn Bay  - Barry Town United: -:-'"
$ws.Cells.Item(18, 2).Value = "Colwyn Bay"
$ws.Cells.Item(18, 3).Value = 42
$ws.Cells.Item(18, 4).Value = 81
$ws.Cells.Item(18, 5).Value = 85
$ws.Cells.Item(18, 6).Value = 1.23
$ws.Cells.Item(18, 7).ClearContents()

# Row 19
$ws.Cells.Item(19, 1).Value = "SSC Napoli  - Como 1907: -:-'"
$ws.Cells.Item(19, 2).Value = "SSC Napoli"
$ws.Cells.Item(19, 3).Value = 42
$ws.Cells.Item(19, 4).Value = 79
$ws.Cells.Item(19, 5).Value = 76
$ws.Cells.Item(19, 6).Value = 1.95
$ws.Cells.Item(19, 7).ClearContents()

# Row 20
$ws.Cells.Item(20, 1).Value = "CD Real Cartagena  - Club Boca Juniors de Cali: 20:05"
$ws.Cells.Item(20, 2).Value = "CD Real Cartagena"
$ws.Cells.Item(20, 3).Value = 41
$ws.Cells.Item(20, 4).Value = 79
$ws.Cells.Item(20, 5).Value = 93
$ws.Cells.Item(20, 6).Value = 1.18
$ws.Cells.Item(20, 7).ClearContents()

# Row 21
$ws.Cells.Item(21, 1).Value = "Al-Khaleej FC - Al-Ittihad Club : 4:0'"
$ws.Cells.Item(21, 2).Value = "Al-Ittihad Club"
$ws.Cells.Item(21, 3).Value = 39
$ws.Cells.Item(21, 4).Value = 86
$ws.Cells.Item(21, 5).Value = 100
$ws.Cells.Item(21, 6).Value = 1.28
$ws.Cells.Item(21, 7).ClearContents()

# Row 22
$ws.Cells.Item(22, 1).Value = "FC Copenhagen  - FC Fredericia: -:-'"
$ws.Cells.Item(22, 2).Value = "FC Copenhagen"
$ws.Cells.Item(22, 3).Value = 38
$ws.Cells.Item(22, 4).Value = 70
$ws.Cells.Item(22, 5).ClearContents()
$ws.Cells.Item(22, 6).Value = 1.18
$ws.Cells.Item(22, 7).ClearContents()

# Row 23
$ws.Cells.Item(23, 1).Value = "Cruzeiro Esporte Clube  - Esporte Clube Vitória: 18:00"
$ws.Cells.Item(23, 2).Value = "Cruzeiro Esporte Clube"
$ws.Cells.Item(23, 3).Value = 37
$ws.Cells.Item(23, 4).ClearContents()
$ws.Cells.Item(23, 5).Value = 96
$ws.Cells.Item(23, 6).Value = 1.44
$ws.Cells.Item(23, 7).ClearContents()

# Row 24
$ws.Cells.Item(24, 1).Value = "Nashville SC - Inter Miami CF : 22:30"
$ws.Cells.Item(24, 2).Value = "Inter Miami CF"
$ws.Cells.Item(24, 3).Value = 36
$ws.Cells.Item(24, 4).Value = 100
$ws.Cells.Item(24, 5).ClearContents()
$ws.Cells.Item(24, 6).Value = 1.44
$ws.Cells.Item(24, 7).ClearContents()

# Row 25
$ws.Cells.Item(25, 1).Value = "Olympiacos Piraeus B  - GS Ilioupolis: 12:00"
$ws.Cells.Item(25, 2).Value = "Olympiacos Piraeus B"
$ws.Cells.Item(25, 3).Value = 34
$ws.Cells.Item(25, 4).Value = 78
$ws.Cells.Item(25, 5).ClearContents()
$ws.Cells.Item(25, 6).Value = 1.3
$ws.Cells.Item(25, 7).ClearContents()

# Row 26
$ws.Cells.Item(26, 1).Value = "Puebla FC - CD Cruz Azul ✓: 0:3"
$ws.Cells.Item(26, 2).Value = "CD Cruz Azul"
$ws.Cells.Item(26, 3).Value = 33
$ws.Cells.Item(26, 4).ClearContents()
$ws.Cells.Item(26, 5).Value = 96
$ws.Cells.Item(26, 6).Value = 1.44
$ws.Cells.Item(26, 7).Value = "✓"

# Row 27
$ws.Cells.Item(27, 1).Value = "CA Cerro - Club Nacional : 21:00"
$ws.Cells.Item(27, 2).Value = "Club Nacional"
$ws.Cells.Item(27, 3).Value = 29
$ws.Cells.Item(27, 4).Value = 86
$ws.Cells.Item(27, 5).Value = 100
$ws.Cells.Item(27, 6).Value = 1.28
$ws.Cells.Item(27, 7).ClearContents()

# Row 28
$ws.Cells.Item(28, 1).Value = "Real Oruro  - Gualberto Villarroel San José: 18:00"
$ws.Cells.Item(28, 2).Value = "Real Oruro"
$ws.Cells.Item(28, 3).Value = 28
$ws.Cells.Item(28, 4).Value = 79
$ws.Cells.Item(28, 5).Value = 93
$ws.Cells.Item(28, 6).Value = 1.18
$ws.Cells.Item(28, 7).ClearContents()
